$d = $word.ActiveDocument

# --- Remove the existing _GoBack bookmark from the heading paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Fix "backwards" -> "backward" in the intro paragraph ---
$d.Content.Find.Execute("In order to provide backwards compatibility", $true, $false, $false, $false, $false, $true, 1, $false, "In order to provide backward compatibility", 2)

# --- Re-insert the _GoBack bookmark right before "backward compatibility..." ---
$r = $d.Content
$r.Find.Execute("backward compatibility with the Speedway Custom ESO version")
$bmRange = $d.Range($r.Start, $r.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Insert a new paragraph with the "Note" text after the intro paragraph ---
$para = $d.Paragraphs(2)
$newPara = $para.Range.InsertParagraphAfter()
$d.Paragraphs(3).Range.Text = "Note:  This process is only applicable when using an “empty” JDA database.  Do not execute the scripts more than once."

# --- "Scripts written" -> "Scripts" heading ---
$d.Content.Find.Execute("Scripts written", $true, $false, $false, $false, $false, $true, 1, $false, "Scripts", 2)

# --- Merge "Creating the " + "POS Options" runs into a single run ---
$d.Content.Find.Execute("Creating the POS Options", $true, $false, $false, $false, $false, $true, 1, $false, "Creating the POS Options", 2)
